$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-Rows($r1, $r2) {
    foreach ($c in $cols) {
        $addr1 = "$c$r1"
        $addr2 = "$c$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Betting odds for these matches were re-scraped/updated; the rows for the
# same fixtures got reshuffled in the source data pull (pairwise swaps).
Swap-Rows 49 50
Swap-Rows 51 52
Swap-Rows 71 73
Swap-Rows 72 74

# Append the newly-scraped fixture as row 100, copying formats from the
# last existing data row (99) so number formats / borders match.
$ws.Range("A99:V99").Copy()
$ws.Range("A100:V100").PasteSpecial(-4122)

$ws.Range("A100").Value = 99
$ws.Range("B100").Value = "wales"
$ws.Range("C100").Value = "cymru-premier"
$ws.Range("D100").Value = "2023-2024"
$ws.Range("E100").Value = 45259.86458333334
$ws.Range("F100").Value = "Colwyn Bay"
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = "Connahs Q."
$ws.Range("I100").Value = 3
$ws.Range("J100").Value = 6.42
$ws.Range("K100").Value = "27/11/2023 21:12"
$ws.Range("L100").Value = 8.82
$ws.Range("M100").Value = "29/11/2023 20:43"
$ws.Range("N100").Value = 5.03
$ws.Range("O100").Value = "27/11/2023 21:12"
$ws.Range("P100").Value = 6.21
$ws.Range("Q100").Value = "29/11/2023 20:43"
$ws.Range("R100").Value = 1.34
$ws.Range("S100").Value = "27/11/2023 21:12"
$ws.Range("T100").Value = 1.27
$ws.Range("U100").Value = "29/11/2023 20:43"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/wales/cymru-premier/colwyn-bay-connahs-q/vPQxW5U7/"
